$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data (and its
# formatting) from A:D to B:E.
$ws.Columns("A:A").Insert()

# New column A header/value for the "tab name" pair.
$ws.Range("A1").Value = "TabName"
$ws.Range("A2").Value = "CasesTab"

# Updated Neo4j "query" text (column B, row 2) reflecting the new Cypher.
$ws.Range("B2").Value = "MATCH (ct:clinical_trial)<--(a:arm)<--(c:case)
    WHERE c.ethnicity IN ['HISPANIC_OR_LATINO']
WITH DISTINCT c, a, ct
RETURN 
    COALESCE(c.case_id, '') AS ``Case ID``,
    COALESCE(ct.clinical_trial_designation, '') AS ``Trial Code``,
    COALESCE(a.arm_id, '') AS ``Arm``,
    COALESCE(a.arm_drug, '') AS ``Arm Treatment``,
    COALESCE(c.disease, '') AS ``Diagnosis``,
    COALESCE(c.gender, '') AS ``Gender``,
    COALESCE(c.race, '') AS ``Race``,
    COALESCE(c.ethnicity, '') AS ``Ethnicity``"

# Updated Neo4j "StatQuery" text (column C, row 2).
$ws.Range("C2").Value = "MATCH (s:specimen)-->(c:case)-->(:arm)-->(ct:clinical_trial)
    WHERE WHERE c.ethnicity IN ['HISPANIC_OR_LATINO']
OPTIONAL MATCH (f:file)-->(:sequencing_assay)-->(:nucleic_acid)-->(s)
RETURN 
    COUNT(DISTINCT f) AS number_of_files,
    COUNT(DISTINCT c.case_id) AS number_of_cases,
    COUNT(DISTINCT ct.clinical_trial_designation) AS number_of_trials"

# New column A gets a narrow best-fit width; the other columns keep the
# widths they already had (shifted right one slot by the Insert above).
$ws.Columns("A:A").ColumnWidth = 8

# Row 2 grows taller to fit the longer wrapped query text.
$ws.Rows(2).RowHeight = 174
